# Gamma ammo/weapon table: add the OG-7B rocket row (row 44) with tripled
# NPC damage and an increased price, matching the author's commit:
#   "Tripled rockets damage to NPCs, increased price to 12k"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- New row 44: ammo_og-7b / EX -------------------------------------------
$ws.Range("A44").Value = "ammo_og-7b"
$ws.Range("B44").Value = "EX"
$ws.Range("C44").Value = 12490
$ws.Range("D44").Formula = "=C44/30"
$ws.Range("E44").Formula = "=K44/D44"
$ws.Range("F44").Formula = "=G44/D44*100"
$ws.Range("G44").Value = 0.37
$ws.Range("H44").Value = 3
$ws.Range("I44").Value = 3
$ws.Range("J44").Formula = "=I44*H44"
$ws.Range("K44").Formula = "=J44*Feuil2!`$B`$1"

# --- Carry down the number formats used by the row above (D/E/F/K), same
# --- as Excel's own "extend formatting" behaviour when continuing a table.
$ws.Range("D43").Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("E43").Copy()
$ws.Range("E44").PasteSpecial(-4122)
$ws.Range("F43").Copy()
$ws.Range("F44").PasteSpecial(-4122)
$ws.Range("K43").Copy()
$ws.Range("K44").PasteSpecial(-4122)

# The PasteSpecial above only touches formatting, but re-assert the formulas
# to be safe regardless of paste semantics.
$ws.Range("D44").Formula = "=C44/30"
$ws.Range("E44").Formula = "=K44/D44"
$ws.Range("F44").Formula = "=G44/D44*100"
$ws.Range("K44").Formula = "=J44*Feuil2!`$B`$1"

# --- Match the selection left behind in the saved view ---------------------
[void]$ws.Range("D26").Select()
